# Update the Entsoe "Notified Production Solar" model for Horeco:
# shift the timestamp column (A) forward by one day and replace the
# notified-production values (B) with the refreshed figures for rows 2:97.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 96,2
$arr[0,0] = 46074.01041666666
$arr[0,1] = 0
$arr[1,0] = 46074.02083333334
$arr[1,1] = 0
$arr[2,0] = 46074.03125
$arr[2,1] = 0
$arr[3,0] = 46074.04166666666
$arr[3,1] = 0
$arr[4,0] = 46074.05208333334
$arr[4,1] = 0.39
$arr[5,0] = 46074.0625
$arr[5,1] = 0
$arr[6,0] = 46074.07291666666
$arr[6,1] = 0
$arr[7,0] = 46074.08333333334
$arr[7,1] = 0
$arr[8,0] = 46074.09375
$arr[8,1] = 0
$arr[9,0] = 46074.10416666666
$arr[9,1] = 0
$arr[10,0] = 46074.11458333334
$arr[10,1] = 0
$arr[11,0] = 46074.125
$arr[11,1] = 0
$arr[12,0] = 46074.13541666666
$arr[12,1] = 0.55
$arr[13,0] = 46074.14583333334
$arr[13,1] = 0
$arr[14,0] = 46074.15625
$arr[14,1] = 0
$arr[15,0] = 46074.16666666666
$arr[15,1] = 0
$arr[16,0] = 46074.17708333334
$arr[16,1] = 0
$arr[17,0] = 46074.1875
$arr[17,1] = 0
$arr[18,0] = 46074.19791666666
$arr[18,1] = 0
$arr[19,0] = 46074.20833333334
$arr[19,1] = 0
$arr[20,0] = 46074.21875
$arr[20,1] = 0.575
$arr[21,0] = 46074.22916666666
$arr[21,1] = 0.573
$arr[22,0] = 46074.23958333334
$arr[22,1] = 0.578
$arr[23,0] = 46074.25
$arr[23,1] = 0.594
$arr[24,0] = 46074.26041666666
$arr[24,1] = 1.683
$arr[25,0] = 46074.27083333334
$arr[25,1] = 4.246
$arr[26,0] = 46074.28125
$arr[26,1] = 9.119999999999999
$arr[27,0] = 46074.29166666666
$arr[27,1] = 15.828
$arr[28,0] = 46074.30208333334
$arr[28,1] = 31.607
$arr[29,0] = 46074.3125
$arr[29,1] = 45.919
$arr[30,0] = 46074.32291666666
$arr[30,1] = 61.422
$arr[31,0] = 46074.33333333334
$arr[31,1] = 77.72199999999999
$arr[32,0] = 46074.34375
$arr[32,1] = 105.236
$arr[33,0] = 46074.35416666666
$arr[33,1] = 127.13
$arr[34,0] = 46074.36458333334
$arr[34,1] = 146.141
$arr[35,0] = 46074.375
$arr[35,1] = 167.626
$arr[36,0] = 46074.38541666666
$arr[36,1] = 191.936
$arr[37,0] = 46074.39583333334
$arr[37,1] = 210.212
$arr[38,0] = 46074.40625
$arr[38,1] = 229.776
$arr[39,0] = 46074.41666666666
$arr[39,1] = 248.431
$arr[40,0] = 46074.42708333334
$arr[40,1] = 265.764
$arr[41,0] = 46074.4375
$arr[41,1] = 280.047
$arr[42,0] = 46074.44791666666
$arr[42,1] = 293.174
$arr[43,0] = 46074.45833333334
$arr[43,1] = 309.361
$arr[44,0] = 46074.46875
$arr[44,1] = 317.62
$arr[45,0] = 46074.47916666666
$arr[45,1] = 324.257
$arr[46,0] = 46074.48958333334
$arr[46,1] = 325.086
$arr[47,0] = 46074.5
$arr[47,1] = 326.904
$arr[48,0] = 46074.51041666666
$arr[48,1] = 320.176
$arr[49,0] = 46074.52083333334
$arr[49,1] = 314.317
$arr[50,0] = 46074.53125
$arr[50,1] = 307.116
$arr[51,0] = 46074.54166666666
$arr[51,1] = 299.66
$arr[52,0] = 46074.55208333334
$arr[52,1] = 284.634
$arr[53,0] = 46074.5625
$arr[53,1] = 272.046
$arr[54,0] = 46074.57291666666
$arr[54,1] = 257.026
$arr[55,0] = 46074.58333333334
$arr[55,1] = 238.529
$arr[56,0] = 46074.59375
$arr[56,1] = 211.995
$arr[57,0] = 46074.60416666666
$arr[57,1] = 190.132
$arr[58,0] = 46074.61458333334
$arr[58,1] = 164.753
$arr[59,0] = 46074.625
$arr[59,1] = 142.254
$arr[60,0] = 46074.63541666666
$arr[60,1] = 114.653
$arr[61,0] = 46074.64583333334
$arr[61,1] = 91.95099999999999
$arr[62,0] = 46074.65625
$arr[62,1] = 71.45399999999999
$arr[63,0] = 46074.66666666666
$arr[63,1] = 53.452
$arr[64,0] = 46074.67708333334
$arr[64,1] = 36.273
$arr[65,0] = 46074.6875
$arr[65,1] = 23.779
$arr[66,0] = 46074.69791666666
$arr[66,1] = 15.575
$arr[67,0] = 46074.70833333334
$arr[67,1] = 9.452999999999999
$arr[68,0] = 46074.71875
$arr[68,1] = 1.313
$arr[69,0] = 46074.72916666666
$arr[69,1] = 1.091
$arr[70,0] = 46074.73958333334
$arr[70,1] = 1.055
$arr[71,0] = 46074.75
$arr[71,1] = 0
$arr[72,0] = 46074.76041666666
$arr[72,1] = 0.55
$arr[73,0] = 46074.77083333334
$arr[73,1] = 0
$arr[74,0] = 46074.78125
$arr[74,1] = 0
$arr[75,0] = 46074.79166666666
$arr[75,1] = 0
$arr[76,0] = 46074.80208333334
$arr[76,1] = 0.39
$arr[77,0] = 46074.8125
$arr[77,1] = 0
$arr[78,0] = 46074.82291666666
$arr[78,1] = 0
$arr[79,0] = 46074.83333333334
$arr[79,1] = 0
$arr[80,0] = 46074.84375
$arr[80,1] = 0.55
$arr[81,0] = 46074.85416666666
$arr[81,1] = 0
$arr[82,0] = 46074.86458333334
$arr[82,1] = 0
$arr[83,0] = 46074.875
$arr[83,1] = 0
$arr[84,0] = 46074.88541666666
$arr[84,1] = 0
$arr[85,0] = 46074.89583333334
$arr[85,1] = 0
$arr[86,0] = 46074.90625
$arr[86,1] = 0
$arr[87,0] = 46074.91666666666
$arr[87,1] = 0
$arr[88,0] = 46074.92708333334
$arr[88,1] = 0
$arr[89,0] = 46074.9375
$arr[89,1] = 0
$arr[90,0] = 46074.94791666666
$arr[90,1] = 0
$arr[91,0] = 46074.95833333334
$arr[91,1] = 0
$arr[92,0] = 46074.96875
$arr[92,1] = 0
$arr[93,0] = 46074.97916666666
$arr[93,1] = 0
$arr[94,0] = 46074.98958333334
$arr[94,1] = 0
$arr[95,0] = 46075
$arr[95,1] = 0

$ws.Range("A2:B97").Value = $arr
